# Update the LR-pairs (Rspo3-Lgr4) sheet with new TPM-derived values.
#
# The "MuSCs" sending-cluster block (previously rows 10-13) is replaced by a
# "Resolving-Mac" sending-cluster block, and a new "Inflammatory-Mac" target
# cluster is inserted for every sending cluster (ECs, FAPs, Resolving-Mac),
# growing each 4-row block to 5 rows (12 data rows -> 15 data rows, i.e.
# rows 2-16 instead of 2-13). All of the numeric NATMI metrics are
# recomputed against the new TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param(
        [int]$Row,
        [object[]]$Values
    )
    for ($i = 0; $i -lt $Values.Count; $i++) {
        $ws.Cells.Item($Row, $i + 1).Value = $Values[$i]
    }
}

# Columns: A Sending cluster | B Ligand symbol | C Receptor symbol | D Target cluster
# E Ligand-expressing cells | F Ligand detection rate | G Ligand average expression value
# H Ligand total expression value | I Ligand derived specificity of average expression value
# J Ligand derived specificity of total expression value | K Receptor-expressing cells
# L Receptor detection rate | M Receptor average expression value | N Receptor total expression value
# O Receptor derived specificity of average expression value | P Receptor derived specificity of total expression value
# Q Edge average expression weight | R Edge total expression weight
# S Edge average expression derived specificity | T Edge total expression derived specificity

Set-RowValues 2 @("ECs", "Rspo3", "Lgr4", "ECs", 1.0, 0.3333333333333333, 0.06762866666666667, 0.202886, 0.0134153952845566, 0.0134153952845566, 3.0, 1.0, 0.5779736666666667, 1.733921, 0.04835019606981441, 0.05356228267519154, 0.03908758844511111, 0.351788296006, 0.0006486369923623752, 0.0007185591944308523)
Set-RowValues 3 @("ECs", "Rspo3", "Lgr4", "FAPs", 1.0, 0.3333333333333333, 0.06762866666666667, 0.202886, 0.0134153952845566, 0.0134153952845566, 3.0, 1.0, 7.791016, 23.373048, 0.6517548686181108, 0.7220131735856595, 0.5268960240586666, 4.742064216528, 0.00874354919114621, 0.009686092124308801)
Set-RowValues 4 @("ECs", "Rspo3", "Lgr4", "Inflammatory-Mac", 1.0, 0.3333333333333333, 0.06762866666666667, 0.202886, 0.0134153952845566, 0.0134153952845566, 1.0, 0.3333333333333333, 0.05015166666666667, 0.150455, 0.004195421100317676, 0.004647681895481942, 0.003391690347777778, 0.03052521313, 0.00005628323244593102, 0.00006235048978476752)
Set-RowValues 5 @("ECs", "Rspo3", "Lgr4", "MuSCs", 1.0, 0.3333333333333333, 0.06762866666666667, 0.202886, 0.0134153952845566, 0.0134153952845566, 2.0, 1.0, 3.489664, 6.979328, 0.2919266886169084, 0.215597330685123, 0.2360013234346667, 1.416007940608, 0.003916311921907495, 0.002892323413436189)
Set-RowValues 6 @("ECs", "Rspo3", "Lgr4", "Resolving-Mac", 1.0, 0.3333333333333333, 0.06762866666666667, 0.202886, 0.0134153952845566, 0.0134153952845566, 1.0, 0.3333333333333333, 0.0451, 0.1353, 0.003772825594848836, 0.004179531158543795, 0.003050052866666667, 0.0274504758, 0.00005061394669458952, 0.00005607006259598581)
Set-RowValues 7 @("FAPs", "Rspo3", "Lgr4", "ECs", 3.0, 1.0, 4.954393, 14.863179, 0.9827953701592058, 0.9827953701592059, 3.0, 1.0, 0.5779736666666667, 1.733921, 0.04835019606981441, 0.05356228267519154, 2.863508688317667, 25.771578194859, 0.04751834884370343, 0.05264076342833689)
Set-RowValues 8 @("FAPs", "Rspo3", "Lgr4", "FAPs", 3.0, 1.0, 4.954393, 14.863179, 0.9827953701592058, 0.9827953701592059, 3.0, 1.0, 7.791016, 23.373048, 0.6517548686181108, 0.7220131735856595, 38.599755133288, 347.397796199592, 0.6405416673566008, 0.7095912041939413)
Set-RowValues 9 @("FAPs", "Rspo3", "Lgr4", "Inflammatory-Mac", 3.0, 1.0, 4.954393, 14.863179, 0.9827953701592058, 0.9827953701592059, 1.0, 0.3333333333333333, 0.05015166666666667, 0.150455, 0.004195421100317676, 0.004647681895481942, 0.2484710662716667, 2.236239596445, 0.004123240433260454, 0.004567720248852415)
Set-RowValues 10 @("FAPs", "Rspo3", "Lgr4", "MuSCs", 3.0, 1.0, 4.954393, 14.863179, 0.9827953701592058, 0.9827953701592059, 2.0, 1.0, 3.489664, 6.979328, 0.2919266886169084, 0.215597330685123, 17.289166893952, 103.735001363712, 0.2869041979986057, 0.2118880584160222)
Set-RowValues 11 @("FAPs", "Rspo3", "Lgr4", "Resolving-Mac", 3.0, 1.0, 4.954393, 14.863179, 0.9827953701592058, 0.9827953701592059, 1.0, 0.3333333333333333, 0.0451, 0.1353, 0.003772825594848836, 0.004179531158543795, 0.2234431243, 2.0109881187, 0.003707915527035587, 0.004107623872052984)
Set-RowValues 12 @("Resolving-Mac", "Rspo3", "Lgr4", "ECs", 1.0, 0.3333333333333333, 0.019102, 0.057306, 0.003789234556237495, 0.003789234556237496, 3.0, 1.0, 0.5779736666666667, 1.733921, 0.04835019606981441, 0.05356228267519154, 0.01104045298066667, 0.09936407682600001, 0.0001832102337485991, 0.0002029600524237967)
Set-RowValues 13 @("Resolving-Mac", "Rspo3", "Lgr4", "FAPs", 1.0, 0.3333333333333333, 0.019102, 0.057306, 0.003789234556237495, 0.003789234556237496, 3.0, 1.0, 7.791016, 23.373048, 0.6517548686181108, 0.7220131735856595, 0.148823987632, 1.339415888688, 0.002469652070363774, 0.002735877267409482)
Set-RowValues 14 @("Resolving-Mac", "Rspo3", "Lgr4", "Inflammatory-Mac", 1.0, 0.3333333333333333, 0.019102, 0.057306, 0.003789234556237495, 0.003789234556237496, 1.0, 0.3333333333333333, 0.05015166666666667, 0.150455, 0.004195421100317676, 0.004647681895481942, 0.0009579971366666668, 0.008621974230000001, 0.00001589743461129167, 0.00001761115684475956)
Set-RowValues 15 @("Resolving-Mac", "Rspo3", "Lgr4", "MuSCs", 1.0, 0.3333333333333333, 0.019102, 0.057306, 0.003789234556237495, 0.003789234556237496, 2.0, 1.0, 3.489664, 6.979328, 0.2919266886169084, 0.215597330685123, 0.066659561728, 0.399957370368, 0.001106178696395172, 0.0008169488556646308)
Set-RowValues 16 @("Resolving-Mac", "Rspo3", "Lgr4", "Resolving-Mac", 1.0, 0.3333333333333333, 0.019102, 0.057306, 0.003789234556237495, 0.003789234556237496, 1.0, 0.3333333333333333, 0.0451, 0.1353, 0.003772825594848836, 0.004179531158543795, 0.0008615002, 0.007753501800000001, 0.00001429612111865849, 0.00001583722389482548)
